$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: law-why -------------------------------------------------------
$ws.Range("A26").Value = "law-why"
$why = @"
เพื่อเป็นเครื่องมือส่งเสริมความเจริญเติบโตทางเศรษฐกิจ
 -เพื่อเป็นเครื่องมือในการควบคุมการบริโภคของประชาชนและบำรุงสาธารณูปโภคเละบริการสาธารณะ
 -เพื่อเป็นเครื่องมือในการกระจายรายได้แก่ให้ประชาชนและเป็นการรักษาเสถียรภาพในทางเศรษฐกิจของประเทศ
"@
$ws.Range("B26").Value = $why

# --- Row 27: law-calculate --------------------------------------------------
$ws.Range("A27").Value = "law-calculate"
$ws.Range("A27").Style = "Normal"
$calc = @"
การคำนวณภาษีของบุคคลธรรมดา เงินได้สุทธิซึ่งเป็นฐานภาษีสำหรับคำนวณภาษีเงินได้บุคคลธรรมดาซึ่งมาจากเงินได้พึงประเมินที่หักค่าใช้จ่ายและค่าลดหย่อนเรียบร้อยแล้ว (ค่าใช้จ่าย-ค่าลดหย่อน=เงินได้สุทธิ)
"@
$ws.Range("B27").Value = $calc

# --- Row 28: law-time --------------------------------------------------------
$ws.Range("A28").Value = "law-time"
$ws.Range("A28").WrapText = $true
$time = @"
กฎหมายกำหนดให้บุคคลต้องทำการยื่นเสียภาษีในช่วง 1 มกราคม - 31 มีนาคม ของทุกปี
"@
$ws.Range("B28").Value = $time

# --- Row 29/30: law-salary and law-place -----------------------------------
# Insert in "salary" then "place" order first so the new shared strings land
# at the same indices as the source workbook, then fix up which row shows
# which text (place ends up on row 29, salary ends up on row 30).
$ws.Range("A29").Value = "law-salary"
$salary = @"
หากมีเงินเดือนหรือมีรายได้จากหลายทางเกิน 10,000 บาท/เดือน (120,000 บาท/ปี) ต้องยื่นภาษีทุกคน
"@
$ws.Range("B29").Value = $salary

$ws.Range("A30").Value = "law-place"
$place = @"
1.สำนักงานสรรพากรทุกสาขาทุกเเห่ง 
2.ไปรษณีย์ เเบบลงทะเบียน 
3.ช่องทางออนไลน์ ผ่านเว็บไซต์ของกรมสรรพากร
"@
$ws.Range("B30").Value = $place

# Now swap: row 29 should display "law-place"/place text, row 30 should
# display "law-salary"/salary text (reusing the strings added above).
$ws.Range("A29").Value = "law-place"
$ws.Range("B29").Value = $place
$ws.Range("A30").Value = "law-salary"
$ws.Range("B30").Value = $salary
$ws.Range("A29").Style = "Normal"

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(26).RowHeight = 102
$ws.Rows.Item(27).RowHeight = 63.75
$ws.Rows.Item(28).RowHeight = 25.5
$ws.Rows.Item(29).RowHeight = 38.25
$ws.Rows.Item(30).RowHeight = 38.25

# --- Sheet view: scroll position + selection --------------------------------
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("B30").Select()
